$p = $ppt.ActivePresentation

# --- Slide 3: "Prior Work" ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange
$tr3.Text = "Many of the existing fracture detection models use some kind of deep learning model to correctly detect fractures`rOf the ones that do not their problem’s can vary from the speed of the algorithm to the accuracy`rOne method utilizing a hough transform is limited to fractures located close to the middle of the bone only`rOf the non deep learning based methods, our method is able to find fractures in a bone at any location"

# --- Slide 4: "Dataset" ---
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(4, 1)
$para4.InsertBefore("Bones featuring no abnormality are placed in folders labelled negative and those with abnormalities are placed in folders labelled positive`r")

# --- Slide 7: "Remaining Work" ---
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(2)
$tr7 = $shape7.TextFrame.TextRange
$tr7.Text = "Extend this program by making use of deep learning`rCreate a more complex algorithm that can detect fractures on bones that have existing metal braces/screws`rExtend the programs capcity to effectively evaluate fractures of any type of bone"
